# modificando template ppt e capa
#
# 1) Update the cached "datetimeFigureOut" footer-date text (23/03/2021 ->
#    03/05/2021) on every slide layout and on the slide master.
# 2) Re-label the six "divider" rectangles on slide 1 and set the two
#    cover textboxes ("titulo" placeholders) to "Bolsa" / "Atleta".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder (ppPlaceholderDate = 16) text on slide master and
#    every custom (slide) layout.
# ---------------------------------------------------------------------
$newDate = "03/05/2021"

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Type -eq 14) {
        if ($shp.PlaceholderFormat.Type -eq 16) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Type -eq 14) {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Slide 1 content.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

# "Retangulo 3" : CONTEMPLADOS -> GERAL (bold, DIN)
$shp = $s.Shapes.Item(1)
$shp.TextFrame.TextRange.Text = "GERAL"
$shp.TextFrame.TextRange.Font.Bold = $true
$shp.TextFrame.TextRange.Font.Name = "DIN"

# "Retangulo 4" : (empty) -> BENEFICIÁRIOS (bold, DIN)
$shp = $s.Shapes.Item(2)
$shp.TextFrame.TextRange.Text = "BENEFICIÁRIOS"
$shp.TextFrame.TextRange.LanguageID = "pt-BR"
$shp.TextFrame.TextRange.Font.Bold = $true
$shp.TextFrame.TextRange.Font.Name = "DIN"

# "Retangulo 5" : (empty) -> ENTIDADES (bold)
$shp = $s.Shapes.Item(3)
$shp.TextFrame.TextRange.Text = "ENTIDADES"
$shp.TextFrame.TextRange.LanguageID = "pt-BR"
$shp.TextFrame.TextRange.Font.Bold = $true

# "Retangulo 6" : (empty) -> QUESTIONÁRIO (bold)
$shp = $s.Shapes.Item(4)
$shp.TextFrame.TextRange.Text = "QUESTIONÁRIO"
$shp.TextFrame.TextRange.LanguageID = "pt-BR"
$shp.TextFrame.TextRange.Font.Bold = $true

# "Retangulo 7" : (empty) -> DADOS - INSCRIÇÕES (bold)
$shp = $s.Shapes.Item(5)
$shp.TextFrame.TextRange.LanguageID = "pt-BR"
$shp.TextFrame.TextRange.Font.Bold = $true
$shp.TextFrame.TextRange.Text = "DADOS - INSCRIÇÕES"

# "Retangulo 8" : (empty) -> DADOS - QUESTIONÁRIO (not bold)
$shp = $s.Shapes.Item(6)
$shp.TextFrame.TextRange.Text = "DADOS - QUESTIONÁRIO"
$shp.TextFrame.TextRange.LanguageID = "pt-BR"

# "CaixaDeTexto 12" : título -> Bolsa
$shp = $s.Shapes.Item(8)
$shp.TextFrame.TextRange.Text = "Bolsa"

# "CaixaDeTexto 14" : título -> Atleta
$shp = $s.Shapes.Item(9)
$shp.TextFrame.TextRange.Text = "Atleta"
